$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update LCD module description and part numbers (row 4)
$ws.Range("A4").Value = "2004 I2C LCD"
$ws.Range("C4").Value = "4411-CN0296D-ND"
$ws.Range("D4").Value = "CN0296D"

# Update note about pullup resistor (row 8, column E)
$ws.Range("E8").Value = "Optional for pullup I2C"

# Set column E width to match new content
$ws.Columns.Item(5).ColumnWidth = 17.63
